# Update Column C (Reaction_number) values for rows 2-20 on both the
# "NBR" and "BAR" sheets, reflecting the new ecoli study / NB ratio analysis
# results described in the commit message.

$wb = $excel.ActiveWorkbook

$wsNBR = $wb.Worksheets.Item("NBR")
$wsBAR = $wb.Worksheets.Item("BAR")

$nbrValues = @(839, 836, 826, 822, 800, 795, 788, 785, 781, 780, 535, 531, 530, 527, 527, 525, 525, 527, 527)
$barValues = @(746, 753, 743, 744, 745, 741, 744, 737, 739, 739, 735, 734, 732, 737, 740, 742, 742, 738, 738)

for ($i = 0; $i -lt $nbrValues.Length; $i++) {
    $row = $i + 2
    $wsNBR.Cells.Item($row, 3).Value = $nbrValues[$i]
    $wsBAR.Cells.Item($row, 3).Value = $barValues[$i]
}
